$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6026200866783623
$ws.Range("B3").Value = 0.5798319321640421
$ws.Range("B4").Value = 0.8269927534359164
$ws.Range("B5").Value = 0.5910059235362268
$ws.Range("B6").Value = 0.5312215685844421
$ws.Range("B7").Value = 0.3766456842422485
